$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "IPP241-Solar_3-Wind_1"
$ws.Range("B2").Value = 129.7316043227053
$ws.Range("C2").Value = 88.24017042041037
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 3164.081687917208
$ws.Range("F2").Value = 4164.081687917208
$ws.Range("G2").Value = 65.00000000011914
$ws.Range("H2").Value = 17.64700479759088

# Row 3
$ws.Range("A3").Value = "IPP241-Solar_3-Wind_6"
$ws.Range("B3").Value = 133.1287247831341
$ws.Range("C3").Value = 87.03541073781425
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3220.761226087495
$ws.Range("F3").Value = 4220.761226087495
$ws.Range("G3").Value = 65.00000000011914
$ws.Range("H3").Value = 18.25556568411022

# Row 4
$ws.Range("A4").Value = "IPP241-Solar_4-Wind_1"
$ws.Range("B4").Value = 103.099190244392
$ws.Range("C4").Value = 99.64452975237268
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 3389.959170067466
$ws.Range("F4").Value = 4389.959170067466
$ws.Range("G4").Value = 65.00000000011939
$ws.Range("H4").Value = 13.50218545543646

# Row 5
$ws.Range("A5").Value = "IPP241-Solar_4-Wind_6"
$ws.Range("B5").Value = 106.191019047453
$ws.Range("C5").Value = 98.05969899680531
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 3453.868883793133
$ws.Range("F5").Value = 4453.868883793133
$ws.Range("G5").Value = 65.0000000001194
$ws.Range("H5").Value = 13.87958373221947
